$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data (rows 2-4), shifting existing data down
$ws.Rows("2:4").Insert()

# Fill the newly inserted rows with new values
$ws.Cells.Item(2, 1).Value = -0.7805059552192688
$ws.Cells.Item(2, 2).Value = -0.7913583517074585
$ws.Cells.Item(2, 3).Value = 0.3286340832710266
$ws.Cells.Item(3, 1).Value = 0.1883520781993866
$ws.Cells.Item(3, 2).Value = 0.3149188160896301
$ws.Cells.Item(3, 3).Value = -0.798948347568512
$ws.Cells.Item(4, 1).Value = 2.627274990081787
$ws.Cells.Item(4, 2).Value = -2.136787414550781
$ws.Cells.Item(4, 3).Value = -1.515871286392212

# Append 7 new rows of data at the end (rows 25-31)
$ws.Cells.Item(25, 1).Value = 6.472879886627197
$ws.Cells.Item(25, 2).Value = -9.945176124572754
$ws.Cells.Item(25, 3).Value = 9.769540786743164
$ws.Cells.Item(26, 1).Value = 3.245661020278931
$ws.Cells.Item(26, 2).Value = -6.337657928466797
$ws.Cells.Item(26, 3).Value = 0.1816275864839553
$ws.Cells.Item(27, 1).Value = 5.774599075317383
$ws.Cells.Item(27, 2).Value = 6.254833221435547
$ws.Cells.Item(27, 3).Value = 0.6913566589355469
$ws.Cells.Item(28, 1).Value = -2.965895891189575
$ws.Cells.Item(28, 2).Value = -1.341034770011902
$ws.Cells.Item(28, 3).Value = 1.717206358909607
$ws.Cells.Item(29, 1).Value = -4.298542022705078
$ws.Cells.Item(29, 2).Value = 1.24329674243927
$ws.Cells.Item(29, 3).Value = 0.5667206645011902
$ws.Cells.Item(30, 1).Value = -6.537195205688477
$ws.Cells.Item(30, 2).Value = -4.836699962615967
$ws.Cells.Item(30, 3).Value = -2.582733631134033
$ws.Cells.Item(31, 1).Value = -1.381847739219666
$ws.Cells.Item(31, 2).Value = 14.15004062652588
$ws.Cells.Item(31, 3).Value = 5.833388328552246
